# Add a new LeetCode entry ("Valid Palindrome") as day 28 / row 30 on
# Sheet1, mirroring how the tracker gains one new row per solved problem.
# Sheet2's rolling averages are plain formulas, so they recalc on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- day number --------------------------------------------------------
$ws.Range("A30").Value = 28

# --- problem name, as a real hyperlink ---------------------------------
# Add the hyperlink first (Excel needs somewhere to attach it); give it the
# address as the display/TextToDisplay for now so the <hyperlink display=.../>
# attribute matches Excel's usual "display = target address" convention,
# then overwrite the cell text with the friendly problem title - this
# keeps the hyperlink relationship while showing "Valid Palindrome" in the
# grid, exactly like the other rows above it.
$ws.Hyperlinks.Add($ws.Range("B30"), "https://leetcode.com/problems/valid-palindrome/", "", "", "https://leetcode.com/problems/valid-palindrome/")
$ws.Range("B30").Value = "Valid Palindrome"
$ws.Range("B30").Style = "Hyperlink"

# --- remaining stats for the row ---------------------------------------
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 59
$ws.Range("F30").Value = 0.78
$ws.Range("G30").Value = 16.9
$ws.Range("H30").Value = 0.0533
$ws.Range("I30").Value = "https://leetcode.com/problems/valid-palindrome/submissions/1070463341/"

# --- selection left over from the author's last save --------------------
$ws.Range("H33").Select()
